$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, pushing existing rows 116:125 down to 117:126
$ws.Rows.Item(116).Insert()

# Populate the new row 116 with the new weekly price entry
$ws.Cells.Item(116, 1).Value = 10
$ws.Cells.Item(116, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(116, 3).Value = "La Araucanía"
$ws.Cells.Item(116, 4).Value = 45166
$ws.Cells.Item(116, 5).Value = 9
$ws.Cells.Item(116, 6).Value = "Fruta"
$ws.Cells.Item(116, 7).Value = 100108
$ws.Cells.Item(116, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(116, 9).Value = 100108007
$ws.Cells.Item(116, 10).Value = "Coco"
$ws.Cells.Item(116, 11).Value = "Sin especificar"
$ws.Cells.Item(116, 12).Value = "Primera"
$ws.Cells.Item(116, 13).Value = 15
$ws.Cells.Item(116, 14).Value = 36000
$ws.Cells.Item(116, 15).Value = 36000
$ws.Cells.Item(116, 16).Value = 36000
$ws.Cells.Item(116, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(116, 18).Value = "Perú"
$ws.Cells.Item(116, 19).Value = 1800
$ws.Cells.Item(116, 20).Value = 20
